$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: title (D21) and link (E21)
$ws.Range("D21").Value = "[Python] 파이썬 Thread 사용법 (Feat. Lock )"
$ws.Range("E21").Value = "https://ms-review.tistory.com/25"

# Row 23: title (D23) and link (E23)
$ws.Range("D23").Value = "[iou 추가] How to get accuracy, F1, precision and recall, iou, for a keras model?"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2858"

# Row 44: title (D44) and link (E44)
$ws.Range("D44").Value = "Google AI Blog(논문) 리뷰: SoundStream (An End-to-End Neural Audio Codec)"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/90"
